$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45. This shifts the existing row 45
# (which contains "Accuracy over PyType" in E45 and 0 in F45) down to row 46,
# carrying its values along with it.
$ws.Rows.Item(45).Insert()

# Copy the formatting (fill/style) from the row above so the new row 45
# matches the rest of the table's look (white-fill cell style).
$ws.Range("A44:F44").Copy()
$ws.Range("A45:F45").PasteSpecial(-4122)  # xlPasteFormats

# Populate the newly inserted row 45 with the new "Scalpel Accuracy:" stat.
$ws.Cells.Item(45, 3).Value = "Scalpel Accuracy:"
$ws.Cells.Item(45, 4).Value = 2000
